# Insert two new weekly records (week of serial date 44979) right before the
# existing row 53, shifting all subsequent rows down by two positions.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(53).Insert()
$ws.Rows.Item(53).Insert()

# New row 53: "Primera" quality record for the new week.
$ws.Cells.Item(53,1).Value  = 7
$ws.Cells.Item(53,2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(53,3).Value  = "Ñuble"
$ws.Cells.Item(53,4).Value  = 44979
$ws.Cells.Item(53,5).Value  = 16
$ws.Cells.Item(53,6).Value  = 100112040
$ws.Cells.Item(53,7).Value  = "Cilantro"
$ws.Cells.Item(53,8).Value  = "Sin especificar"
$ws.Cells.Item(53,9).Value  = "Primera"
$ws.Cells.Item(53,10).Value = 300
$ws.Cells.Item(53,11).Value = 700
$ws.Cells.Item(53,12).Value = 800
$ws.Cells.Item(53,13).Value = 750
$ws.Cells.Item(53,14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(53,15).Value = "Provincia de Diguillín"
$ws.Cells.Item(53,16).Value = 750
$ws.Cells.Item(53,17).Value = 1
$ws.Cells.Item(53,18).Value = "Hortaliza"

# New row 54: "Segunda" quality record for the new week.
$ws.Cells.Item(54,1).Value  = 7
$ws.Cells.Item(54,2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(54,3).Value  = "Ñuble"
$ws.Cells.Item(54,4).Value  = 44979
$ws.Cells.Item(54,5).Value  = 16
$ws.Cells.Item(54,6).Value  = 100112040
$ws.Cells.Item(54,7).Value  = "Cilantro"
$ws.Cells.Item(54,8).Value  = "Sin especificar"
$ws.Cells.Item(54,9).Value  = "Segunda"
$ws.Cells.Item(54,10).Value = 200
$ws.Cells.Item(54,11).Value = 600
$ws.Cells.Item(54,12).Value = 600
$ws.Cells.Item(54,13).Value = 600
$ws.Cells.Item(54,14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(54,15).Value = "Provincia de Diguillín"
$ws.Cells.Item(54,16).Value = 600
$ws.Cells.Item(54,17).Value = 1
$ws.Cells.Item(54,18).Value = "Hortaliza"
